# Auto-generated edit script applying numeric corrections to multiple
# worksheets per the commit diff ("chore: update Sheets via scheduled runner").
# For each affected row, columns H..N (currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ), LeveProfit(NQ/HQ)) are corrected to their new values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1500
$ws.Range("I31").Value = 1500
$ws.Range("K31").Value = 4500
$ws.Range("M31").Value = -4270
$ws.Range("H98").Value = 6807295
$ws.Range("I98").Value = 9525569
$ws.Range("K98").Value = 9525569
$ws.Range("M98").Value = -9524071
$ws.Range("H122").Value = 6807295
$ws.Range("I122").Value = 9525569
$ws.Range("K122").Value = 28576707
$ws.Range("M122").Value = -28574257
$ws.Range("H123").Value = 98999.664
$ws.Range("J123").Value = 98999.664
$ws.Range("L123").Value = 98999.664
$ws.Range("N123").Value = -108799.664
$ws.Range("H133").Value = 93497.5
$ws.Range("J133").Value = 93497.5
$ws.Range("L133").Value = 93497.5
$ws.Range("N133").Value = -103617.5
$ws.Range("H135").Value = 4657.375
$ws.Range("I135").Value = 1106.6666
$ws.Range("J135").Value = 15309.5
$ws.Range("K135").Value = 9959.999400000001
$ws.Range("L135").Value = 137785.5
$ws.Range("M135").Value = -7424.999400000001
$ws.Range("N135").Value = -142855.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13925.179
$ws.Range("I32").Value = 15414.417
$ws.Range("K32").Value = 15414.417
$ws.Range("M32").Value = -15127.417
$ws.Range("H74").Value = 1973.3334
$ws.Range("I74").Value = 1631.5
$ws.Range("K74").Value = 1631.5
$ws.Range("M74").Value = -757.5
$ws.Range("H77").Value = 1973.3334
$ws.Range("I77").Value = 1631.5
$ws.Range("K77").Value = 8157.5
$ws.Range("M77").Value = -3789.5
$ws.Range("H122").Value = 5355.778
$ws.Range("I122").Value = 4868.1665
$ws.Range("J122").Value = 6331
$ws.Range("K122").Value = 14604.4995
$ws.Range("L122").Value = 18993
$ws.Range("M122").Value = -12154.4995
$ws.Range("N122").Value = -23893
$ws.Range("H124").Value = 17856.5
$ws.Range("J124").Value = 17856.5
$ws.Range("L124").Value = 17856.5
$ws.Range("N124").Value = -27676.5
$ws.Range("H125").Value = 48398.6
$ws.Range("J125").Value = 48398.6
$ws.Range("L125").Value = 48398.6
$ws.Range("N125").Value = -58238.6
$ws.Range("H127").Value = 88247.5
$ws.Range("J127").Value = 88247.5
$ws.Range("L127").Value = 88247.5
$ws.Range("N127").Value = -98167.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 5322.769
$ws.Range("I11").Value = 4320.1
$ws.Range("K11").Value = 4320.1
$ws.Range("M11").Value = -4180.1
$ws.Range("H68").Value = 71330
$ws.Range("J68").Value = 71330
$ws.Range("L68").Value = 71330
$ws.Range("N68").Value = -72952
$ws.Range("H71").Value = 71330
$ws.Range("J71").Value = 71330
$ws.Range("L71").Value = 213990
$ws.Range("N71").Value = -222102
$ws.Range("H76").Value = 46000
$ws.Range("J76").Value = 46000
$ws.Range("L76").Value = 46000
$ws.Range("N76").Value = -46630
$ws.Range("H79").Value = 46000
$ws.Range("J79").Value = 46000
$ws.Range("L79").Value = 46000
$ws.Range("N79").Value = -48184
$ws.Range("H86").Value = 9335.3125
$ws.Range("J86").Value = 2525.7144
$ws.Range("L86").Value = 2525.7144
$ws.Range("N86").Value = -4771.7144
$ws.Range("H89").Value = 9335.3125
$ws.Range("J89").Value = 2525.7144
$ws.Range("L89").Value = 12628.572
$ws.Range("N89").Value = -23860.572
$ws.Range("H124").Value = 78689.5
$ws.Range("J124").Value = 78689.5
$ws.Range("L124").Value = 78689.5
$ws.Range("N124").Value = -88509.5
$ws.Range("H125").Value = 99999.664
$ws.Range("J125").Value = 99999.664
$ws.Range("L125").Value = 99999.664
$ws.Range("N125").Value = -109839.664
$ws.Range("H126").Value = 79999.75
$ws.Range("J126").Value = 79999.75
$ws.Range("L126").Value = 79999.75
$ws.Range("N126").Value = -89879.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 439.85715
$ws.Range("I7").Value = 96.5
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 96.5
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 16.5
$ws.Range("N7").Value = -2726
$ws.Range("H99").Value = 2870.8096
$ws.Range("I99").Value = 2839.35
$ws.Range("K99").Value = 2839.35
$ws.Range("M99").Value = -1341.35
$ws.Range("H122").Value = 4287.6665
$ws.Range("I122").Value = 4231
$ws.Range("J122").Value = 4486
$ws.Range("K122").Value = 12693
$ws.Range("L122").Value = 13458
$ws.Range("M122").Value = -10243
$ws.Range("N122").Value = -18358
$ws.Range("H126").Value = 2870.8096
$ws.Range("I126").Value = 2839.35
$ws.Range("K126").Value = 8518.049999999999
$ws.Range("M126").Value = -6048.049999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12994337
$ws.Range("J4").Value = 7050284.5
$ws.Range("L4").Value = 21150853.5
$ws.Range("N4").Value = -21151077.5
$ws.Range("H6").Value = 2499
$ws.Range("I6").Value = 2499
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 7497
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -7384
$ws.Range("N6").ClearContents()
$ws.Range("H80").Value = 15000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 45000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -46872
$ws.Range("H83").Value = 15000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 135000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -144360
$ws.Range("H124").Value = 11203.5
$ws.Range("I124").Value = 12346.333
$ws.Range("J124").Value = 7775
$ws.Range("K124").Value = 37038.999
$ws.Range("L124").Value = 23325
$ws.Range("M124").Value = -32128.999
$ws.Range("N124").Value = -33145
$ws.Range("H125").Value = 4095.125
$ws.Range("I125").Value = 3072.6
$ws.Range("J125").Value = 5799.3335
$ws.Range("K125").Value = 9217.799999999999
$ws.Range("L125").Value = 17398.0005
$ws.Range("M125").Value = -4297.799999999999
$ws.Range("N125").Value = -27238.0005
$ws.Range("H126").Value = 7269.3335
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 9904
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 29712
$ws.Range("M126").Value = -1060
$ws.Range("N126").Value = -39592
$ws.Range("H129").Value = 1179.9474
$ws.Range("I129").Value = 447.23077
$ws.Range("K129").Value = 1341.69231
$ws.Range("M129").Value = 3658.30769
$ws.Range("H130").Value = 9721.25
$ws.Range("I130").Value = 8999.5
$ws.Range("J130").Value = 10443
$ws.Range("K130").Value = 26998.5
$ws.Range("L130").Value = 31329
$ws.Range("M130").Value = -21978.5
$ws.Range("N130").Value = -41369
$ws.Range("H131").Value = 1946
$ws.Range("J131").Value = 1946
$ws.Range("L131").Value = 5838
$ws.Range("N131").Value = -15918

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 100000
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H95").Value = 48999.5
$ws.Range("J95").Value = 48999.5
$ws.Range("L95").Value = 48999.5
$ws.Range("N95").Value = -54491.5
$ws.Range("H97").Value = 4100.25
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H102").Value = 50866.383
$ws.Range("I102").Value = 3475
$ws.Range("J102").Value = 114054.89
$ws.Range("K102").Value = 3475
$ws.Range("L102").Value = 114054.89
$ws.Range("M102").Value = -1853
$ws.Range("N102").Value = -117298.89
$ws.Range("H122").Value = 2834.8235
$ws.Range("I122").Value = 2689.1428
$ws.Range("K122").Value = 8067.428400000001
$ws.Range("M122").Value = -5617.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2418.5652
$ws.Range("I22").Value = 891.25
$ws.Range("K22").Value = 891.25
$ws.Range("M22").Value = -596.25
$ws.Range("H27").Value = 2418.5652
$ws.Range("I27").Value = 891.25
$ws.Range("K27").Value = 891.25
$ws.Range("M27").Value = -784.25
$ws.Range("H40").Value = 4665.7646
$ws.Range("I40").Value = 4434.6665
$ws.Range("K40").Value = 4434.6665
$ws.Range("M40").Value = -4298.6665
$ws.Range("H122").Value = 3653.1667
$ws.Range("I122").Value = 2204.4443
$ws.Range("K122").Value = 6613.3329
$ws.Range("M122").Value = -4163.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2688903.2
$ws.Range("I113").Value = 4902288
$ws.Range("K113").Value = 14706864
$ws.Range("M113").Value = -14704694
$ws.Range("H136").Value = 13294.353
$ws.Range("J136").Value = 1000
$ws.Range("L136").Value = 3000
$ws.Range("N136").Value = -8100
$ws.Range("H140").Value = 99713
$ws.Range("J140").Value = 99713
$ws.Range("L140").Value = 99713
$ws.Range("N140").Value = -110073
